$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits at the end
#    of the "JS30 - Episode 8, 2:21" paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Insert a brand-new list paragraph right after the "Practice "
#    heading, re-using the numbering/list formatting (numId 1) that
#    is already used by the other bullet items in the document. The
#    easiest reliable way to clone that formatting is to copy an
#    existing ListParagraph paragraph (including its end-of-paragraph
#    mark) and paste it into place, then replace its text.
# ------------------------------------------------------------------

# Locate the "Practice " heading paragraph.
$practiceIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Practice") {
        $practiceIndex = $i
        break
    }
}

# Locate the paragraph immediately after it (an existing ListParagraph
# bullet) so we can clone its paragraph formatting / numbering.
$templateIndex = $practiceIndex + 1
$templatePara = $d.Paragraphs($templateIndex)
$templateRange = $d.Range($templatePara.Range.Start, $templatePara.Range.End)
$templateRange.Copy()

# Paste the cloned paragraph right after "Practice ".
$practicePara = $d.Paragraphs($practiceIndex)
$insertPoint = $d.Range($practicePara.Range.End, $practicePara.Range.End)
$insertPoint.Paste()

# The pasted paragraph is now the new paragraph at $templateIndex.
$newIndex = $templateIndex
$newPara = $d.Paragraphs($newIndex)

# Replace its text (keep the paragraph mark) with the new sentence,
# appending one temporary placeholder character so that the real
# "end of sentence" position is not the literal last character of the
# paragraph (inserting/bookmarking exactly at that spot is unreliable).
$newText = "Hold shift to select multiple checkboxes " + [char]0x2013 + " episode 10"
$bodyRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$bodyRange.Text = $newText + "X"

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark, collapsed, at the end of the
#    new sentence (i.e. just before the placeholder "X").
# ------------------------------------------------------------------
$newPara2 = $d.Paragraphs($newIndex)
$bookmarkPos = $newPara2.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary placeholder character now that the bookmark is
# anchored safely at the desired (non-boundary) position.
$newPara3 = $d.Paragraphs($newIndex)
$placeholderRange = $d.Range($newPara3.Range.End - 2, $newPara3.Range.End - 1)
$placeholderRange.Delete()
